$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Agt"
$ws.Cells.Item(2,3).Value = "Agtr1a"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.4683593333333333
$ws.Cells.Item(2,8).Value = 1.405078
$ws.Cells.Item(2,9).Value = 0.6051469521021553
$ws.Cells.Item(2,10).Value = 0.6051469521021552
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.7246186666666666
$ws.Cells.Item(2,14).Value = 2.173856
$ws.Cells.Item(2,15).Value = 0.02655111241446272
$ws.Cells.Item(2,16).Value = 0.02655111241446271
$ws.Cells.Item(2,17).Value = 0.3393819156408889
$ws.Cells.Item(2,18).Value = 3.054437240768
$ws.Cells.Item(2,19).Value = 0.01606732475253381
$ws.Cells.Item(2,20).Value = 0.0160673247525338

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Agt"
$ws.Cells.Item(3,3).Value = "Agtr1a"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.4683593333333333
$ws.Cells.Item(3,8).Value = 1.405078
$ws.Cells.Item(3,9).Value = 0.6051469521021553
$ws.Cells.Item(3,10).Value = 0.6051469521021552
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 17.195945
$ws.Cells.Item(3,14).Value = 51.587835
$ws.Cells.Item(3,15).Value = 0.630085160334334
$ws.Cells.Item(3,16).Value = 0.630085160334334
$ws.Cells.Item(3,17).Value = 8.053881336236666
$ws.Cells.Item(3,18).Value = 72.48493202613
$ws.Cells.Item(3,19).Value = 0.3812941143411201
$ws.Cells.Item(3,20).Value = 0.38129411434112

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Agt"
$ws.Cells.Item(4,3).Value = "Agtr1a"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.4683593333333333
$ws.Cells.Item(4,8).Value = 1.405078
$ws.Cells.Item(4,9).Value = 0.6051469521021553
$ws.Cells.Item(4,10).Value = 0.6051469521021552
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 9.370898
$ws.Cells.Item(4,14).Value = 28.112694
$ws.Cells.Item(4,15).Value = 0.3433637272512032
$ws.Cells.Item(4,16).Value = 0.3433637272512032
$ws.Cells.Item(4,17).Value = 4.388947540014667
$ws.Cells.Item(4,18).Value = 39.500527860132
$ws.Cells.Item(4,19).Value = 0.2077855130085014
$ws.Cells.Item(4,20).Value = 0.2077855130085013

# Row 5
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Agt"
$ws.Cells.Item(5,3).Value = "Agtr1a"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.3056003333333333
$ws.Cells.Item(5,8).Value = 0.916801
$ws.Cells.Item(5,9).Value = 0.3948530478978448
$ws.Cells.Item(5,10).Value = 0.3948530478978448
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.7246186666666666
$ws.Cells.Item(5,14).Value = 2.173856
$ws.Cells.Item(5,15).Value = 0.02655111241446272
$ws.Cells.Item(5,16).Value = 0.02655111241446271
$ws.Cells.Item(5,17).Value = 0.2214437060728889
$ws.Cells.Item(5,18).Value = 1.992993354656
$ws.Cells.Item(5,19).Value = 0.01048378766192891
$ws.Cells.Item(5,20).Value = 0.01048378766192891

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Agt"
$ws.Cells.Item(6,3).Value = "Agtr1a"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.3056003333333333
$ws.Cells.Item(6,8).Value = 0.916801
$ws.Cells.Item(6,9).Value = 0.3948530478978448
$ws.Cells.Item(6,10).Value = 0.3948530478978448
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 17.195945
$ws.Cells.Item(6,14).Value = 51.587835
$ws.Cells.Item(6,15).Value = 0.630085160334334
$ws.Cells.Item(6,16).Value = 0.630085160334334
$ws.Cells.Item(6,17).Value = 5.255086523981666
$ws.Cells.Item(6,18).Value = 47.295778715835
$ws.Cells.Item(6,19).Value = 0.248791045993214
$ws.Cells.Item(6,20).Value = 0.248791045993214

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Agt"
$ws.Cells.Item(7,3).Value = "Agtr1a"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.3056003333333333
$ws.Cells.Item(7,8).Value = 0.916801
$ws.Cells.Item(7,9).Value = 0.3948530478978448
$ws.Cells.Item(7,10).Value = 0.3948530478978448
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 9.370898
$ws.Cells.Item(7,14).Value = 28.112694
$ws.Cells.Item(7,15).Value = 0.3433637272512032
$ws.Cells.Item(7,16).Value = 0.3433637272512032
$ws.Cells.Item(7,17).Value = 2.863749552432667
$ws.Cells.Item(7,18).Value = 25.773745971894
$ws.Cells.Item(7,19).Value = 0.1355782142427019
$ws.Cells.Item(7,20).Value = 0.1355782142427019
